$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.840.66"
$ws.Range("E2").Value = "  -2.11%  "

$ws.Range("D3").Value = "1.806.11"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "308.94"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "0.9978"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("E7").Value = "  +3.66%  "

$ws.Range("E8").Value = "  -0.81%  "

$ws.Range("D9").Value = "0.07300"
$ws.Range("E9").Value = "  -2.93%  "

$ws.Range("D10").Value = "0.8602"
$ws.Range("E10").Value = "  -3.64%  "

$ws.Range("D11").Value = "20.44"

$ws.Range("D12").Value = "1.748.49"
$ws.Range("E12").Value = "  -4.34%  "

$ws.Range("D13").Value = "5.331"
$ws.Range("E13").Value = "  -1.47%  "

$ws.Range("D14").Value = "6.511"
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").Value = "0.07051"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").Value = "91.37"
$ws.Range("E16").Value = "  -3.26%  "

$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").Value = "0.000008674"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "0.9991"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").Value = "14.69"
$ws.Range("E20").Value = "  -3.59%  "

$ws.Range("D21").Value = "26.817.88"
$ws.Range("E21").Value = "  -2.24%  "

$ws.Range("D22").Value = "5.298"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("E23").Value = "  -2.30%  "

$ws.Range("D24").Value = "1.930.98"
$ws.Range("E24").Value = "  -6.04%  "

$ws.Range("D25").Value = "1.908"
$ws.Range("E25").Value = "  -3.77%  "

$ws.Range("D26").Value = "150.79"
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("D27").Value = "18.30"
$ws.Range("E27").Value = "  -1.60%  "

$ws.Range("D28").Value = "2.142"
$ws.Range("E28").Value = "  -10.06%  "

$ws.Range("D29").Value = "5.247"
$ws.Range("E29").Value = "  -2.29%  "

$ws.Range("D30").Value = "114.64"
$ws.Range("E30").Value = "  -2.53%  "

$ws.Range("D31").Value = "0.08908"
$ws.Range("E31").Value = "  +0.83%  "

$ws.Range("D32").Value = "0.7655"
$ws.Range("E32").Value = "  -2.29%  "

$ws.Range("E33").Value = "  -3.14%  "

$ws.Range("D34").Value = "4.455"
$ws.Range("E34").Value = "  -1.65%  "

$ws.Range("D35").Value = "2.888"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "0.9985"
$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("D37").Value = "1.120"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").Value = "0.01947"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "2.430"
$ws.Range("E39").Value = "  +5.86%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.05222"
$ws.Range("E40").Value = "  -2.05%  "

$ws.Range("D41").Value = "2.899"
$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("D42").Value = "7.182"
$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("D43").Value = "0.5233"
$ws.Range("E43").Value = "  -1.56%  "

$ws.Range("D44").Value = "0.1657"
$ws.Range("E44").Value = "  -4.32%  "

$ws.Range("D45").Value = "8.533"
$ws.Range("E45").Value = "  -3.00%  "

$ws.Range("D46").Value = "0.5039"
$ws.Range("E46").Value = "  -1.78%  "

$ws.Range("D47").Value = "10.33"
$ws.Range("E47").Value = "  -3.53%  "

$ws.Range("D48").Value = "104.31"
$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("D49").Value = "0.9983"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").Value = "1.654"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("D51").Value = "0.06301"
$ws.Range("E51").Value = "  -1.17%  "
